# many many bug fixes and a reasonable set of lifemap specs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural edits -------------------------------------------------
# New input column "in:total_order_price" goes in before the old column D
# (out:Shipping Method), shifting everything from D: onward one column right.
$ws.Columns("D:D").EntireColumn.Insert()

# A new pricing-tier rule row is added right after row 4 (between the two
# "FedEx International Express Styrofoam Box" rows), shifting rows 5: down.
$ws.Rows("5:5").EntireRow.Insert()

# The new column is narrower than the "out:Shipping Method" column that used
# to live there.
$ws.Columns("D:D").ColumnWidth = 17.2

# --- header row ---------------------------------------------------------
$ws.Range("D2").Value = 'in:total_order_price'

# --- rows 7 & 8 (accumulate rule, no longer conditioned on quantity) -------
$ws.Range("B7").Value = '$(Recombinant protein) intersection $in #=0'
$ws.Range("B8").Value = '$(Recombinant protein) intersection $in #=0'

# --- row 3 ---------------------------------------------------------------
$ws.Range("B3").Value = '$in >= $(Recombinant protein)'
$ws.Range("D3").Value = '<2500.01'

# --- row 4 ---------------------------------------------------------------
$ws.Range("B4").Value = '$in >= $(Recombinant protein)'
$ws.Range("C4").Value = '<7'
$ws.Range("D4").Value = '>2500.00'

# --- row 5 (newly inserted) ----------------------------------------------
$ws.Range("A5").Value = '$(ProSpec)'
$ws.Range("B5").Value = '$in >= $(Recombinant protein)'
$ws.Range("C5").Value = '>6'
$ws.Range("E5").Value = 'FedEx International Express Styrofoam Box'
$ws.Range("F5").Value = 'ice packs/blue ice'
$ws.Range("G5").Value = 2000
$ws.Range("H5").Value = 12000
$ws.Range("I5").Value = 13000
$ws.Range("J5").Value = 17500
$ws.Range("K5").Value = 17500
$ws.Range("L5").Value = 25000

# --- row 6 -----------------------------------------------------------------
$ws.Range("B6").Value = '$in >= $(Recombinant protein)'
$ws.Range("C6").Value = '<7'
$ws.Range("D6").Value = '<2500.01'
$ws.Range("I6").Value = 'N/A'
$ws.Range("J6").Value = 'N/A'
$ws.Range("K6").Value = 'N/A'
$ws.Range("L6").Value = 'N/A'

# --- leave the selection where the author left it -------------------------
$ws.Range("M6").Select()
